$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Rename sheet "REGCV2" -> "REGCV1"
# ---------------------------------------------------------------------------
$wsRegcv = $wb.Worksheets.Item("REGCV2")
$wsRegcv.Name = "REGCV1"

# ---------------------------------------------------------------------------
# 2) REGCV1 sheet: insert four new parameter columns (Kpid, Kiid, Kpiq, Kiiq)
#    before the existing Kpvd.. columns, shifting O:V -> S:Z, and update the
#    values in row 2 to match the new layout.
# ---------------------------------------------------------------------------

# Header row (row 1), columns O..Z
$headers = @("Kpid", "Kiid", "Kpiq", "Kiiq", "Kpvd", "Kivd", "Kpvq", "Kivq", "ra", "xs", "gammap", "gammaq")
$startCol = 15  # column O
for ($i = 0; $i -lt $headers.Count; $i++) {
    $cell = $wsRegcv.Cells.Item(1, $startCol + $i)
    $cell.NumberFormat = "@"
    $cell.Value = $headers[$i]
}

# Data row (row 2), columns O..Z
$values = @("100", "200", "100", "200", "1000", "500", "1000", "500", "0", "0.15", "1", "1")
for ($i = 0; $i -lt $values.Count; $i++) {
    $cell = $wsRegcv.Cells.Item(2, $startCol + $i)
    $cell.NumberFormat = "@"
    $cell.Value = $values[$i]
}

# ---------------------------------------------------------------------------
# 3) Toggler sheet: append a new row 4 (second toggler tied to GENROU_2)
# ---------------------------------------------------------------------------
$wsToggler = $wb.Worksheets.Item("Toggler")

$row4 = @("2", "3", "1", "Toggler_2", "SynGen", "GENROU_2", "1")
for ($i = 0; $i -lt $row4.Count; $i++) {
    $cell = $wsToggler.Cells.Item(4, $i + 1)
    $cell.NumberFormat = "@"
    $cell.Value = $row4[$i]
}
